# Apply updated cryptocurrency price/volume data to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.675.37'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.505.32'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.24'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.57'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.503.39'
$ws.Range('E9').Value = '  -0.64%  '
$ws.Range('E10').Value = '  -1.20%  '
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.358'
$ws.Range('E12').Value = '  +2.83%  '
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.961.39'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '69.533.92'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('E16').Value = '  +0.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.73'
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.511.12'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.22'
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('E20').Value = '  -3.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '349.24'
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('E23').Value = '  -1.12%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.86'
$ws.Range('E25').Value = '  +2.19%  '
$ws.Range('E26').Value = '  -1.94%  '
$ws.Range('E27').Value = '  -2.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.630.61'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.44%  '
$ws.Range('E30').Value = '  -1.91%  '
$ws.Range('E31').Value = '  -0.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '458.18'
$ws.Range('E32').Value = '  -1.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.21'
$ws.Range('E33').Value = '  -5.74%  '
$ws.Range('E34').Value = '  -1.39%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '157.52'
$ws.Range('E36').Value = '  +2.72%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.116'
$ws.Range('E37').Value = '  -2.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.05'
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  -1.28%  '
$ws.Range('E42').Value = '  -2.00%  '
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('E45').Value = '  -5.24%  '
$ws.Range('E46').Value = '  -7.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '141.34'
$ws.Range('E47').Value = '  -1.27%  '
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('E49').Value = '  -2.44%  '
$ws.Range('E50').Value = '  -0.67%  '
$ws.Range('E51').Value = '  -0.42%  '
